# "yeni kelime eklendi 28 ekim" - add a new vocabulary word to the list.
# The previous word entry on row 11 (interior / iç, dahili / sıfat) is
# cleared and a new word entry is written one row below, on row 12:
#   intensity | yoğunluk | isim

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old word row entirely.
$ws.Range("C11:E11").ClearContents()

# Write the new word in the next row.
$ws.Range("C12").Value = "intensity"
$ws.Range("D12").Value = "yoğunluk"
$ws.Range("E12").Value = "isim"

# Move the active selection to the newly entered word, like a user would
# leave it after typing the entry.
$ws.Range("C12").Select()
